# Leetcode Progress.xlsx -- "Aton picked 3 linked list problems"
# Rebuild the single "Sheet1" into three sheets: Easy / Medium / Hard, with
# Easy carrying the updated LeetCode tracking data (3 new rows assigned to
# Aton using Python), and Medium / Hard as blank templates that mirror the
# Easy sheet's header + row formatting.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the existing sheet to "Easy" and clear it out so we can
#    rewrite its contents from a clean slate (keeps the pre-existing
#    column widths).
# ---------------------------------------------------------------------
$easy = $wb.Worksheets.Item(1)
$easy.Name = "Easy"
$easy.UsedRange.Clear()

# ---------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------
$headers = @("#", "Question", "Topic", "PIC", "Status", "Language", "Level", "Comment")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $easy.Cells.Item(1, $col)
    $cell.Value2 = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108  # xlCenter
}

# ---------------------------------------------------------------------
# 3. Data rows: #, Question, Topic, PIC, Status, Language, Level, Comment
# ---------------------------------------------------------------------
$rows = @(
    @(21,  "Merge Two Sort Lists",                "Linked List", "Pedoe", "Ongoing", "",       "Easy", ""),
    @(83,  "Remove Duplicates from Sorted List",   "Linked List", "Aton",  "Ongoing", "Python", "Easy", ""),
    @(141, "Linked List Cycle",                    "Linked List", "Pedoe", "Ongoing", "",       "Easy", "Solve it without using extra space"),
    @(206, "Reverse Linked List",                  "Linked List", "Pedoe", "Ongoing", "",       "Easy", "Reverse a single linked list"),
    @(234, "Palindrone Linked List",                "Linked List", "Aton",  "Ongoing", "Python", "Easy", "O(n) time complexity, O(1) space complexity"),
    @(237, "Delete Node in a Linked List",          "Linked List", "Aton",  "Ongoing", "Python", "Easy", "")
)

# Seed new shared strings in the same order the original author's edit
# would have (Status text, then the "Python" language tag, then the new
# question titles by row order) before filling in the rest of the grid.
$easy.Cells.Item(3, 5).Value2 = "Ongoing"
$easy.Cells.Item(3, 6).Value2 = "Python"
$easy.Cells.Item(7, 2).Value2 = "Delete Node in a Linked List"
$easy.Cells.Item(3, 2).Value2 = "Remove Duplicates from Sorted List"

$r = 2
foreach ($row in $rows) {
    $easy.Cells.Item($r, 1).Value2 = $row[0]
    $easy.Cells.Item($r, 2).Value2 = $row[1]
    $easy.Cells.Item($r, 3).Value2 = $row[2]
    $easy.Cells.Item($r, 4).Value2 = $row[3]

    $statusCell = $easy.Cells.Item($r, 5)
    $statusCell.Value2 = $row[4]
    $statusCell.Font.Bold = $true
    $statusCell.Font.ThemeColor = 6   # -> theme "5" (Accent2)

    $langCell = $easy.Cells.Item($r, 6)
    if ($row[5] -ne "") {
        $langCell.Value2 = $row[5]
        $langCell.Font.Bold = $true
        $langCell.Font.ThemeColor = 5   # -> theme "4" (Accent1)
    } else {
        $langCell.Font.Bold = $true
    }

    $easy.Cells.Item($r, 7).Value2 = $row[6]

    if ($row[7] -ne "") {
        $easy.Cells.Item($r, 8).Value2 = $row[7]
    }

    $r++
}

# Rows for "Palindrone Linked List" / "Delete Node in a Linked List" (Aton,
# col B) keep the original left-aligned Question style.
$easy.Range("B6").HorizontalAlignment = -4131  # xlLeft
$easy.Range("B7").HorizontalAlignment = -4131  # xlLeft

$easy.Range("E10").Select()

# ---------------------------------------------------------------------
# 4. Medium / Hard: blank templates that mirror Easy's header + row
#    formatting, with no data filled in yet.
# ---------------------------------------------------------------------
$medium = $wb.Worksheets.Add($null, $easy)
$medium.Name = "Medium"
$hard = $wb.Worksheets.Add($null, $medium)
$hard.Name = "Hard"

foreach ($tpl in @($medium, $hard)) {
    for ($col = 1; $col -le $headers.Length; $col++) {
        $cell = $tpl.Cells.Item(1, $col)
        $cell.Value2 = $headers[$col - 1]
        $cell.Font.Bold = $true
        $cell.HorizontalAlignment = -4108  # xlCenter
    }

    $tpl.Range("B2").HorizontalAlignment = -4131
    $tpl.Range("B3").HorizontalAlignment = -4131

    foreach ($rr in @(2, 3, 5, 6, 7)) {
        $statusCell = $tpl.Cells.Item($rr, 5)
        $statusCell.Font.Bold = $true
        $statusCell.Font.ThemeColor = 6
    }

    foreach ($rr in @(2, 3)) {
        $langCell = $tpl.Cells.Item($rr, 6)
        $langCell.Font.Bold = $true
        $langCell.Font.ThemeColor = 5
    }
    foreach ($rr in @(5, 6, 7)) {
        $langCell = $tpl.Cells.Item($rr, 6)
        $langCell.Font.Bold = $true
    }

    $tpl.Range("D11").Select()
}

$easy.Select()
